$d = $word.ActiveDocument

$pairs = @(
    @("42×30=", "60×21="),
    @("19×35=", "47×35="),
    @("35×58=", "91×57="),
    @("39×73=", "25×46="),
    @("58×37=", "37×73="),
    @("90×38=", "12×22="),
    @("38×80=", "97×92="),
    @("95×56=", "41×52="),
    @("74×22=", "85×75="),
    @("78×94=", "85×84="),
    @("90×33=", "81×22="),
    @("94×41=", "20×70="),
    @("72×81=", "29×55="),
    @("81×14=", "21×43="),
    @("54×91=", "35×43="),
    @("31×69=", "87×35="),
    @("57×22=", "39×16="),
    @("90×97=", "95×80="),
    @("55×67=", "80×64="),
    @("85×29=", "70×61="),
    @("70×60=", "43×14="),
    @("52×76=", "52×54="),
    @("87×25=", "21×84="),
    @("11×16=", "38×61="),
    @("19×82=", "40×47=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
